$wb = $excel.ActiveWorkbook

# Overview sheet has the same shared string text used for B2/C2/B3/C3 -> update via direct value set
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = "Handed back: in sync with en-US"
$ovw.Range("C2").Value = "Handed back: in sync with en-US"
$ovw.Range("B3").Value = "Handed back: in sync with en-US"
$ovw.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("H2").Value = "2016-03-23 11:13:18"
$zh.Range("H3").Value = "2016-03-23 11:13:18"

$zh.Range("F2").Value = "bf25903d-bfde-49bb-a091-e574762c400c.md"
$zh.Range("G2").Value = "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.zh-cn.xlf"
$zh.Range("F3").Value = "bf25903d-bfde-49bb-a091-e574762c400c.md"
$zh.Range("G3").Value = "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.zh-cn.xlf"

$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/df9a81f36d93c6b00e0e05ddcde3c4b1f58e9f1f/e2e/bf25903d-bfde-49bb-a091-e574762c400c.md", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29d39b0b91ae8461a6276685aaafbaf15af49f5a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.zh-cn.xlf", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/df9a81f36d93c6b00e0e05ddcde3c4b1f58e9f1f/e2e/bf25903d-bfde-49bb-a091-e574762c400c.md", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.md")
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/29d39b0b91ae8461a6276685aaafbaf15af49f5a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.zh-cn.xlf", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.zh-cn.xlf")

# de-de sheet
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("H2").Value = "2016-03-23 11:13:25"
$de.Range("H3").Value = "2016-03-23 11:13:25"

$de.Range("F2").Value = "bf25903d-bfde-49bb-a091-e574762c400c.md"
$de.Range("G2").Value = "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.de-de.xlf"
$de.Range("F3").Value = "bf25903d-bfde-49bb-a091-e574762c400c.md"
$de.Range("G3").Value = "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.de-de.xlf"

$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/df9a81f36d93c6b00e0e05ddcde3c4b1f58e9f1f/e2e/bf25903d-bfde-49bb-a091-e574762c400c.md", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74246af896936fb2026db008aca038ebb2289c33/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.de-de.xlf", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/df9a81f36d93c6b00e0e05ddcde3c4b1f58e9f1f/e2e/bf25903d-bfde-49bb-a091-e574762c400c.md", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.md")
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/74246af896936fb2026db008aca038ebb2289c33/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.de-de.xlf", "", "", "bf25903d-bfde-49bb-a091-e574762c400c.314048f49f06a58ff9c768703cf2246ce24a91ca.de-de.xlf")

Write-Host "done"
